$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.582.31'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.959.55'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.42%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.615'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.45%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.42'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +6.07%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +5.52%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0790'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.22%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +7.11%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.840'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.35%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.248.47'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.29'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.03%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.963.02'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.516.44'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.76%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '229.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.78%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.29%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.46'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +6.25%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +5.08%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.145'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +10.73%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.20'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.53'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.18%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +11.10%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.93%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.75'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.18%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.70%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.57'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +21.46%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.28'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +8.30%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.51'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -8.33%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0981'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.84%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.31%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.08%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.72%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.97'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.91%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.369.70'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.31%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.91%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.22%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.14'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.68%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '44.69'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.36%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.139.46'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.41%  '
$ws.Range('E51').Style = 'Normal'
